$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:B22 with the new computed coefficient values
$values = @(
    47.93981441118185,
    0.2707226458805622,
    -0.06386218596864779,
    0.2213498106413643,
    0.3738469168289906,
    0.3632458766623982,
    -0.2032666004689614,
    0.2027482470925074,
    -0.07352680939214168,
    0.1446805381226363,
    0.08546923656109882,
    0.2212358988818773,
    0.05545545341961042,
    -0.1106656443753674,
    -0.05476329283728745,
    0.6196579270378004,
    0.1713491161932482,
    0.0635071712792292,
    -0.1018295869040546,
    -0.3142911166931832,
    -0.2207810511162162
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Remove row 23 entirely (A23/B23), shrinking the used range to A1:B22
$ws.Rows.Item(23).Delete()
